# The sheet contains weekly "Perejil" price records for "Vega Modelo de
# Temuco" ordered (roughly) from most-recent-date to oldest. A new weekly
# record was inserted right before the existing row 397, pushing the old
# rows 397..448 down to 398..449 and growing the used range from
# A1:R448 to A1:R449.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 397; this shifts rows 397-448 down
# to 398-449 and keeps all of their existing values/formatting intact.
$ws.Rows(397).Insert()

# Populate the newly-inserted row 397 with the new weekly record.
$ws.Cells.Item(397, 1).Value = 10
$ws.Cells.Item(397, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(397, 3).Value = "La Araucanía"
$ws.Cells.Item(397, 4).Value = 44984
$ws.Cells.Item(397, 5).Value = 9
$ws.Cells.Item(397, 6).Value = 100112044
$ws.Cells.Item(397, 7).Value = "Perejil"
$ws.Cells.Item(397, 8).Value = "Sin especificar"
$ws.Cells.Item(397, 9).Value = "Primera"
$ws.Cells.Item(397, 10).Value = 50
$ws.Cells.Item(397, 11).Value = 5000
$ws.Cells.Item(397, 12).Value = 5000
$ws.Cells.Item(397, 13).Value = 5000
$ws.Cells.Item(397, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(397, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(397, 16).Value = 1667
$ws.Cells.Item(397, 17).Value = 3
$ws.Cells.Item(397, 18).Value = "Hortaliza"
